# Apply the "end of loop" progress update to the Bill Summary sheet.
# - Columns C8:C17 ("Qty executed upto date") are genuine numbers -> set directly.
# - Columns G/H in rows 9,10,11,13,14,19,21 ("Upto date Amount" / totals) are
#   stored as TEXT that merely looks like a formatted number (e.g. "9472.00").
#   Assigning a numeric-looking string straight to .Value would make Excel
#   silently reinterpret it as a Number (dropping the trailing zero and the
#   decimal formatting), so instead we write it as a string formula
#   (="9472.00") and immediately flatten that formula down to its cached
#   text value with Copy / PasteSpecial(values) - the same trick a user
#   would use in the UI to paste "numbers as text" without leaving a
#   quote-prefix / text-format style behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNumber {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    $cell.Formula = '="' + $Text + '"'
    $cell.Copy()
    $cell.PasteSpecial("values")
}

# --- Quantities executed upto date (plain numbers) ---
$ws.Range("C8").Value = 29
$ws.Range("C9").Value = 37
$ws.Range("C10").Value = 12
$ws.Range("C11").Value = 95
$ws.Range("C12").Value = 34
$ws.Range("C13").Value = 90
$ws.Range("C14").Value = 38
$ws.Range("C15").Value = 84
$ws.Range("C16").Value = 34
$ws.Range("C17").Value = 20

# --- Upto date amounts (text that looks like a formatted number) ---
Set-TextNumber "G9"  "9472.00"
Set-TextNumber "G10" "5664.00"
Set-TextNumber "G11" "62890.00"
Set-TextNumber "G13" "12240.00"
Set-TextNumber "G14" "874.00"

# --- Grand totals (also text) ---
Set-TextNumber "G19" "91140.00"
Set-TextNumber "H19" "91140.00"
Set-TextNumber "G21" "91140.00"
Set-TextNumber "H21" "91140.00"

# Leave Excel's clipboard / marching-ants state clean.
$excel.CutCopyMode = $false
